# Simulated and logged 2021 conference championships
# Update Rams "Players Data" workbook with updated Rushing and Receiving stats.

$wb = $excel.ActiveWorkbook

# ---- Rushing sheet ----
$rushing = $wb.Worksheets.Item("Rushing")

# M.Stafford (row 2): RZATT
$rushing.Range("E2").Value = 12

# S.Michel (row 4): 1DATT, 2DATT, 3DATT, RZATT
$rushing.Range("C4").Value = 108
$rushing.Range("D4").Value = 64
$rushing.Range("E4").Value = 36
$rushing.Range("F4").Value = 39

# C.Akers (row 7): 1DATT, 2DATT, RZATT
$rushing.Range("C7").Value = 40
$rushing.Range("D7").Value = 18
$rushing.Range("F7").Value = 12

# C.Kupp (row 8): 1DATT
$rushing.Range("C8").Value = 3

# ---- Receiving sheet ----
$receiving = $wb.Worksheets.Item("Receiving")

# S.Michel (row 2): Short Target, Short Comp, RZ Target, RZ Comp
$receiving.Range("C2").Value = 51
$receiving.Range("D2").Value = 34
$receiving.Range("G2").Value = 10
$receiving.Range("H2").Value = 4

# C.Akers (row 3): Short Target, Short Comp
$receiving.Range("C3").Value = 7
$receiving.Range("D3").Value = 7

# C.Kupp (row 4): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$receiving.Range("C4").Value = 159
$receiving.Range("D4").Value = 124
$receiving.Range("E4").Value = 52
$receiving.Range("F4").Value = 33
$receiving.Range("G4").Value = 40
$receiving.Range("H4").Value = 28

# V.Jefferson (row 5): Short Target, Short Comp, Deep Target
$receiving.Range("C5").Value = 66
$receiving.Range("D5").Value = 41
$receiving.Range("E5").Value = 31

# B.Skowronek (row 6): Deep Target
$receiving.Range("E6").Value = 5

# O.Beckham (row 7): Short Target, Short Comp, Deep Target, Deep Comp
$receiving.Range("C7").Value = 86
$receiving.Range("D7").Value = 71
$receiving.Range("E7").Value = 34
$receiving.Range("F7").Value = 15

# K.Blanton (row 8): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$receiving.Range("C8").Value = 9
$receiving.Range("D8").Value = 7
$receiving.Range("E8").Value = 1
$receiving.Range("F8").Value = 1
$receiving.Range("G8").Value = 3
$receiving.Range("H8").Value = 2

# T.Higbee (row 9): Short Target, Short Comp, RZ Target, RZ Comp
$receiving.Range("C9").Value = 88
$receiving.Range("D9").Value = 65
$receiving.Range("G9").Value = 20
$receiving.Range("H9").Value = 15
